$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price/Volume columns so numeric-looking strings
# (e.g. "1.00", "0.998") are preserved exactly instead of being coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '67.685.76'
$ws.Range('E2').Value = '  -0.06%  '

$ws.Range('D3').Value = '3.770.79'
$ws.Range('E3').Value = '  -0.79%  '

$ws.Range('D4').Value = '1.01'
$ws.Range('E4').Value = '  +0.50%  '

$ws.Range('D5').Value = '598.87'
$ws.Range('E5').Value = '  +0.35%  '

$ws.Range('D6').Value = '162.82'
$ws.Range('E6').Value = '  -2.66%  '

$ws.Range('D7').Value = '3.767.55'
$ws.Range('E7').Value = '  -0.90%  '

$ws.Range('E8').Value = '  +0.17%  '

$ws.Range('D9').Value = '0.512'
$ws.Range('E9').Value = '  -1.27%  '

$ws.Range('E10').Value = '  -3.16%  '

$ws.Range('E11').Value = '  -1.28%  '

$ws.Range('D12').Value = '6.59'
$ws.Range('E12').Value = '  +4.83%  '

$ws.Range('E13').Value = '  -3.73%  '

$ws.Range('D14').Value = '35.02'
$ws.Range('E14').Value = '  -2.54%  '

$ws.Range('D15').Value = '4.404.06'

$ws.Range('D16').Value = '3.752.58'
$ws.Range('E16').Value = '  -0.92%  '

$ws.Range('D17').Value = '67.746.36'
$ws.Range('E17').Value = '  -0.08%  '

$ws.Range('D18').Value = '18.18'
$ws.Range('E18').Value = '  -1.93%  '

$ws.Range('E19').Value = '  +1.76%  '

$ws.Range('D20').Value = '6.98'
$ws.Range('E20').Value = '  -1.32%  '

$ws.Range('D21').Value = '456.67'
$ws.Range('E21').Value = '  -0.89%  '

$ws.Range('E22').Value = '  -4.36%  '

$ws.Range('D23').Value = '0.690'
$ws.Range('E23').Value = '  -0.48%  '

$ws.Range('D24').Value = '82.66'

$ws.Range('E25').Value = '  -6.47%  '

$ws.Range('D26').Value = '11.81'
$ws.Range('E26').Value = '  -2.07%  '

$ws.Range('E27').Value = '  -0.85%  '

$ws.Range('E28').Value = '  +0.02%  '

$ws.Range('D29').Value = '9.80'
$ws.Range('E29').Value = '  -1.98%  '

$ws.Range('D30').Value = '3.917.68'
$ws.Range('E30').Value = '  -0.75%  '

$ws.Range('E32').Value = '  -2.77%  '

$ws.Range('E33').Value = '  -6.61%  '

$ws.Range('D34').Value = '28.81'
$ws.Range('E34').Value = '  -2.36%  '

$ws.Range('D35').Value = '0.998'
$ws.Range('E35').Value = '  -0.04%  '

$ws.Range('D36').Value = '8.91'
$ws.Range('E36').Value = '  -1.55%  '

$ws.Range('D37').Value = '0.0987'
$ws.Range('E37').Value = '  -1.24%  '

$ws.Range('E38').Value = '  +3.34%  '

$ws.Range('D39').Value = '5.75'
$ws.Range('E39').Value = '  -0.38%  '

$ws.Range('D40').Value = '0.973'
$ws.Range('E40').Value = '  -2.71%  '

$ws.Range('D41').Value = '3.15'
$ws.Range('E41').Value = '  -6.07%  '

$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.10%  '

$ws.Range('D44').Value = '43.43'
$ws.Range('E44').Value = '  +1.71%  '

$ws.Range('D45').Value = '47.13'

$ws.Range('D46').Value = '152.13'
$ws.Range('E46').Value = '  +3.10%  '

$ws.Range('E47').Value = '  -2.00%  '

$ws.Range('D48').Value = '8.26'
$ws.Range('E48').Value = '  -0.85%  '

$ws.Range('E49').Value = '  -0.35%  '

$ws.Range('D50').Value = '1.83'
$ws.Range('E50').Value = '  -0.70%  '

$ws.Range('D51').Value = '385.53'
$ws.Range('E51').Value = '  -2.41%  '
